# Auto-generated Excel COM-interop script to apply the leve-profit data refresh
# (scheduled runner update) described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 3000
$ws.Cells.Item(125, 10).Value = 3000
$ws.Cells.Item(125, 12).Value = 27000
$ws.Cells.Item(125, 14).Value = -31920
$ws.Cells.Item(129, 8).Value = 8914.925999999999
$ws.Cells.Item(129, 9).Value = 409.5
$ws.Cells.Item(129, 10).Value = 12496.158
$ws.Cells.Item(129, 11).Value = 1228.5
$ws.Cells.Item(129, 12).Value = 37488.474
$ws.Cells.Item(129, 13).Value = 3771.5
$ws.Cells.Item(129, 14).Value = -47488.474
$ws.Cells.Item(138, 8).Value = 3127.4753
$ws.Cells.Item(138, 9).Value = 2066.625
$ws.Cells.Item(138, 10).Value = 3504.6667
$ws.Cells.Item(138, 11).Value = 6199.875
$ws.Cells.Item(138, 12).Value = 10514.0001
$ws.Cells.Item(138, 13).Value = -1059.875
$ws.Cells.Item(138, 14).Value = -20794.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1768.909
$ws.Cells.Item(45, 9).Value = 1791.8334
$ws.Cells.Item(45, 10).Value = 1741.4
$ws.Cells.Item(45, 11).Value = 1791.8334
$ws.Cells.Item(45, 12).Value = 1741.4
$ws.Cells.Item(45, 13).Value = -1414.8334
$ws.Cells.Item(45, 14).Value = -2495.4
$ws.Cells.Item(122, 8).Value = 2787.8
$ws.Cells.Item(122, 9).Value = 2166.75
$ws.Cells.Item(122, 10).Value = 5272
$ws.Cells.Item(122, 11).Value = 6500.25
$ws.Cells.Item(122, 12).Value = 15816
$ws.Cells.Item(122, 13).Value = -4050.25
$ws.Cells.Item(122, 14).Value = -20716

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1892.7142
$ws.Cells.Item(99, 9).Value = 1937.5
$ws.Cells.Item(99, 10).Value = 1833
$ws.Cells.Item(99, 11).Value = 1937.5
$ws.Cells.Item(99, 12).Value = 1833
$ws.Cells.Item(99, 13).Value = -439.5
$ws.Cells.Item(99, 14).Value = -4829
$ws.Cells.Item(105, 8).Value = 2128.2144
$ws.Cells.Item(105, 9).Value = 1377.2222
$ws.Cells.Item(105, 11).Value = 1377.2222
$ws.Cells.Item(105, 13).Value = 369.7778000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 39387.06
$ws.Cells.Item(4, 10).Value = 10598.75
$ws.Cells.Item(4, 12).Value = 10598.75
$ws.Cells.Item(4, 14).Value = -10822.75
$ws.Cells.Item(16, 8).Value = 83334664
$ws.Cells.Item(16, 9).Value = 500000000
$ws.Cells.Item(16, 10).Value = 1596
$ws.Cells.Item(16, 11).Value = 500000000
$ws.Cells.Item(16, 12).Value = 1596
$ws.Cells.Item(16, 13).Value = -499999713
$ws.Cells.Item(16, 14).Value = -2170
$ws.Cells.Item(99, 8).Value = 1789226.2
$ws.Cells.Item(99, 9).Value = 3973169.5
$ws.Cells.Item(99, 10).Value = 2363.4546
$ws.Cells.Item(99, 11).Value = 3973169.5
$ws.Cells.Item(99, 12).Value = 2363.4546
$ws.Cells.Item(99, 13).Value = -3971671.5
$ws.Cells.Item(99, 14).Value = -5359.4546
$ws.Cells.Item(105, 8).Value = 962
$ws.Cells.Item(105, 9).Value = 962
$ws.Cells.Item(105, 11).Value = 962
$ws.Cells.Item(105, 13).Value = 785
$ws.Cells.Item(113, 8).Value = 83334664
$ws.Cells.Item(113, 9).Value = 500000000
$ws.Cells.Item(113, 10).Value = 1596
$ws.Cells.Item(113, 11).Value = 500000000
$ws.Cells.Item(113, 12).Value = 1596
$ws.Cells.Item(113, 13).Value = -499997830
$ws.Cells.Item(113, 14).Value = -5936
$ws.Cells.Item(126, 8).Value = 1789226.2
$ws.Cells.Item(126, 9).Value = 3973169.5
$ws.Cells.Item(126, 10).Value = 2363.4546
$ws.Cells.Item(126, 11).Value = 11919508.5
$ws.Cells.Item(126, 12).Value = 7090.3638
$ws.Cells.Item(126, 13).Value = -11917038.5
$ws.Cells.Item(126, 14).Value = -12030.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 802.775
$ws.Cells.Item(5, 10).Value = 1127
$ws.Cells.Item(5, 12).Value = 3381
$ws.Cells.Item(5, 14).Value = -3605
$ws.Cells.Item(107, 8).Value = 464.9
$ws.Cells.Item(107, 10).Value = 629.8
$ws.Cells.Item(107, 12).Value = 1889.4
$ws.Cells.Item(107, 14).Value = -5729.4
$ws.Cells.Item(135, 8).Value = 802.775
$ws.Cells.Item(135, 10).Value = 1127
$ws.Cells.Item(135, 12).Value = 10143
$ws.Cells.Item(135, 14).Value = -15213

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1255.5769
$ws.Cells.Item(102, 9).Value = 1185.8422
$ws.Cells.Item(102, 10).Value = 1444.8572
$ws.Cells.Item(102, 11).Value = 1185.8422
$ws.Cells.Item(102, 12).Value = 1444.8572
$ws.Cells.Item(102, 13).Value = 436.1578
$ws.Cells.Item(102, 14).Value = -4688.8572
$ws.Cells.Item(105, 8).Value = 36500
$ws.Cells.Item(105, 10).Value = 36500
$ws.Cells.Item(105, 12).Value = 36500
$ws.Cells.Item(105, 14).Value = -43488
$ws.Cells.Item(126, 8).Value = 2979.4546
$ws.Cells.Item(126, 9).Value = 1975.3334
$ws.Cells.Item(126, 10).Value = 4184.4
$ws.Cells.Item(126, 11).Value = 5926.0002
$ws.Cells.Item(126, 12).Value = 12553.2
$ws.Cells.Item(126, 13).Value = -3456.0002
$ws.Cells.Item(126, 14).Value = -17493.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1613.0294
$ws.Cells.Item(7, 9).Value = 1313.8
$ws.Cells.Item(7, 10).Value = 1737.7084
$ws.Cells.Item(7, 11).Value = 1313.8
$ws.Cells.Item(7, 12).Value = 1737.7084
$ws.Cells.Item(7, 13).Value = -1201.8
$ws.Cells.Item(7, 14).Value = -1961.7084
$ws.Cells.Item(16, 8).Value = 1400
$ws.Cells.Item(16, 9).Value = 1400
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1400
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -1230
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 12532680
$ws.Cells.Item(68, 9).Value = 33834280
$ws.Cells.Item(68, 10).Value = 2326.8823
$ws.Cells.Item(68, 11).Value = 33834280
$ws.Cells.Item(68, 12).Value = 2326.8823
$ws.Cells.Item(68, 13).Value = -33833531
$ws.Cells.Item(68, 14).Value = -3824.8823
$ws.Cells.Item(71, 8).Value = 12532680
$ws.Cells.Item(71, 9).Value = 33834280
$ws.Cells.Item(71, 10).Value = 2326.8823
$ws.Cells.Item(71, 11).Value = 169171400
$ws.Cells.Item(71, 12).Value = 11634.4115
$ws.Cells.Item(71, 13).Value = -169167656
$ws.Cells.Item(71, 14).Value = -19122.4115
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 1613.0294
$ws.Cells.Item(126, 9).Value = 1313.8
$ws.Cells.Item(126, 10).Value = 1737.7084
$ws.Cells.Item(126, 11).Value = 3941.4
$ws.Cells.Item(126, 12).Value = 5213.1252
$ws.Cells.Item(126, 13).Value = -1471.4
$ws.Cells.Item(126, 14).Value = -10153.1252
$ws.Cells.Item(132, 8).Value = 8777573
$ws.Cells.Item(132, 9).Value = 13164769
$ws.Cells.Item(132, 10).Value = 3180.158
$ws.Cells.Item(132, 11).Value = 39494307
$ws.Cells.Item(132, 12).Value = 9540.474
$ws.Cells.Item(132, 13).Value = -39491777
$ws.Cells.Item(132, 14).Value = -14600.474
$ws.Cells.Item(136, 8).Value = 6183.448
$ws.Cells.Item(136, 9).Value = 11636.454
$ws.Cells.Item(136, 10).Value = 2851.0557
$ws.Cells.Item(136, 11).Value = 34909.362
$ws.Cells.Item(136, 12).Value = 8553.167099999999
$ws.Cells.Item(136, 13).Value = -32359.362
$ws.Cells.Item(136, 14).Value = -13653.1671

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 25471.791
$ws.Cells.Item(100, 9).Value = 54503
$ws.Cells.Item(100, 10).Value = 906.9231
$ws.Cells.Item(100, 11).Value = 109006
$ws.Cells.Item(100, 12).Value = 1813.8462
$ws.Cells.Item(100, 13).Value = -108465
$ws.Cells.Item(100, 14).Value = -2895.8462
$ws.Cells.Item(107, 8).Value = 699.9231
$ws.Cells.Item(107, 9).Value = 283.26315
$ws.Cells.Item(107, 11).Value = 849.78945
$ws.Cells.Item(107, 13).Value = 1070.21055
$ws.Cells.Item(113, 8).Value = 370.63635
$ws.Cells.Item(113, 9).Value = 308.8125
$ws.Cells.Item(113, 10).Value = 535.5
$ws.Cells.Item(113, 11).Value = 926.4375
$ws.Cells.Item(113, 12).Value = 1606.5
$ws.Cells.Item(113, 13).Value = 1243.5625
$ws.Cells.Item(113, 14).Value = -5946.5
$ws.Cells.Item(126, 8).Value = 2992.2
$ws.Cells.Item(126, 9).Value = 3160.2222
$ws.Cells.Item(126, 10).Value = 1480
$ws.Cells.Item(126, 11).Value = 9480.6666
$ws.Cells.Item(126, 12).Value = 4440
$ws.Cells.Item(126, 13).Value = -7010.6666
$ws.Cells.Item(126, 14).Value = -9380
$ws.Cells.Item(136, 8).Value = 307420.47
$ws.Cells.Item(136, 9).Value = 525770.3
$ws.Cells.Item(136, 10).Value = 1730.75
$ws.Cells.Item(136, 11).Value = 1577310.9
$ws.Cells.Item(136, 12).Value = 5192.25
$ws.Cells.Item(136, 13).Value = -1574760.9
$ws.Cells.Item(136, 14).Value = -10292.25

Write-Output "Applied Zeromus_Profits scheduled-runner update."
